# Updated symbol list on Sun Jan 15 23:42:24 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest crypto quotes.
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the workbook's existing inlineStr/text cells)
# instead of auto-converting to numbers or percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.31"
$ws.Range("E2").Value = "'-1.03%"
$ws.Range("D3").Value = "'31.35"
$ws.Range("E3").Value = "'-3.01%"
$ws.Range("D4").Value = "'5.141"
$ws.Range("E4").Value = "'-3.04%"
$ws.Range("D5").Value = "'0.07397"
$ws.Range("E5").Value = "'-2.47%"
$ws.Range("D6").Value = "'2.214"
$ws.Range("E6").Value = "'30.00%"
$ws.Range("D7").Value = "'7.935"
$ws.Range("E7").Value = "'0.63%"
$ws.Range("D8").Value = "'3.827"
$ws.Range("E8").Value = "'-1.01%"
$ws.Range("D9").Value = "'0.9203"
$ws.Range("E9").Value = "'-0.90%"
$ws.Range("D10").Value = "'0.1703"
$ws.Range("E10").Value = "'0.59%"
$ws.Range("D11").Value = "'0.07574"
$ws.Range("E11").Value = "'-4.86%"
$ws.Range("D12").Value = "'0.08134"
$ws.Range("E12").Value = "'1.17%"
$ws.Range("D13").Value = "'0.03018"
$ws.Range("E13").Value = "'-1.37%"
$ws.Range("D14").Value = "'0.09925"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.001493"
$ws.Range("E15").Value = "'-0.69%"
$ws.Range("D16").Value = "'0.006081"
$ws.Range("E16").Value = "'-3.93%"
$ws.Range("E17").Value = "'0.83%"
$ws.Range("E18").Value = "'-0.67%"
$ws.Range("D19").Value = "'0.3262"
$ws.Range("E19").Value = "'-1.18%"
$ws.Range("D20").Value = "'0.1320"
$ws.Range("E20").Value = "'-1.77%"
$ws.Range("D21").Value = "'4.650"
$ws.Range("E21").Value = "'2.16%"
$ws.Range("D22").Value = "'0.04636"
$ws.Range("E22").Value = "'0.94%"
$ws.Range("D23").Value = "'0.1568"
$ws.Range("E23").Value = "'-2.99%"
$ws.Range("D24").Value = "'0.001226"
$ws.Range("E24").Value = "'0.97%"
$ws.Range("D25").Value = "'0.004473"
$ws.Range("E25").Value = "'-0.40%"
$ws.Range("E26").Value = "'-7.05%"
$ws.Range("D27").Value = "'0.0003428"
$ws.Range("E27").Value = "'101.94%"
$ws.Range("D39").Value = "'0.01736"
$ws.Range("E39").Value = "'-0.18%"
$ws.Range("D40").Value = "'0.04510"
$ws.Range("E40").Value = "'-0.39%"
$ws.Range("D41").Value = "'0.007311"
$ws.Range("E41").Value = "'5.08%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("E43").Value = "'7.83%"
$ws.Range("D44").Value = "'0.01062"
$ws.Range("E44").Value = "'-23.69%"
$ws.Range("D45").Value = "'0.00006292"
$ws.Range("E45").Value = "'2.51%"
$ws.Range("E46").Value = "'-22.92%"
$ws.Range("E47").Value = "'12.39%"
